# The sheet is a weekly price-list log. A new weekly record (2022-06-10,
# serial 44722) is inserted above the existing row 55, pushing all the
# rows below it (old 55..70) down by one (new 56..71). The used range
# therefore grows from A1:T70 to A1:T71.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 55; Excel shifts rows 55:70 down to 56:71 and grows
# the sheet dimension automatically.
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with this week's entry. All of the
# "descriptive" columns (market/region/product taxonomy/unit/origin) stay
# constant across this product's rows, only the date + volume/price
# figures change.
$ws.Range("A55").Value = 10
$ws.Range("B55").Value = "Vega Modelo de Temuco"
$ws.Range("C55").Value = "La Araucanía"
$ws.Range("D55").Value = 44722
$ws.Range("E55").Value = 9
$ws.Range("F55").Value = "Fruta"
$ws.Range("G55").Value = 100108
$ws.Range("H55").Value = "Tropicales y subtropicales"
$ws.Range("I55").Value = 100108007
$ws.Range("J55").Value = "Coco"
$ws.Range("K55").Value = "Sin especificar"
$ws.Range("L55").Value = "Primera"
$ws.Range("M55").Value = 15
$ws.Range("N55").Value = 20000
$ws.Range("O55").Value = 20000
$ws.Range("P55").Value = 20000
$ws.Range("Q55").Value = "$/malla 20 unidades"
$ws.Range("R55").Value = "Perú"
$ws.Range("S55").Value = 1000
$ws.Range("T55").Value = 20
